$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tasks")

# Set the cells that introduce brand-new text first, and in this exact
# order, so the rebuilt shared-strings table lines up with the target:
#   ... Background Layout, Levels Layout, Splash Screen Button update,
#       Method: modArray, Accepted
$ws.Range("B16").Value = "Splash Screen Button update"
$ws.Range("B13").Value = "Method: modArray"
$ws.Range("J9").Value = "Accepted"

# --- Row 10: Testing Level 1 marked Completed + Completion Date added ---
$ws.Range("C9").Copy()
$ws.Range("C10").PasteSpecial(-4122)
$ws.Range("C10").Value = "Completed"

$ws.Range("F10").Copy()
$ws.Range("I10").PasteSpecial(-4122)
$ws.Range("I10").Value = 42633

# --- Row 13: marked Completed, dates added ---
$ws.Range("C9").Copy()
$ws.Range("C13").PasteSpecial(-4122)
$ws.Range("C13").Value = "Completed"

$ws.Range("F13").Value = 42631
$ws.Range("G13").Value = 42636

$ws.Range("F10").Copy()
$ws.Range("I13").PasteSpecial(-4122)
$ws.Range("I13").Value = 42633

# --- Row 14: Background Layout, dates added (status/owner unchanged) ---
$ws.Range("F14").Value = 42631
$ws.Range("G14").Value = 42636

# --- Row 15: Levels Layout, marked Completed, dates added ---
$ws.Range("C9").Copy()
$ws.Range("C15").PasteSpecial(-4122)
$ws.Range("C15").Value = "Completed"

$ws.Range("F15").Value = 42631
$ws.Range("G15").Value = 42636

$ws.Range("F10").Copy()
$ws.Range("I15").PasteSpecial(-4122)
$ws.Range("I15").Value = 42633

# --- Row 16: new task "Splash Screen Button update" assigned to Ammar, dates added ---
$ws.Range("C16").Value = "Assigned"
$ws.Range("D16").Value = "Ammar"
$ws.Range("F16").Value = 42631
$ws.Range("G16").Value = 42636

# --- Sheet view / column width tweaks ---
$ws.Columns("B").ColumnWidth = 21.92
$ws.Activate()
$ws.Range("I11").Select()
